# Append the "mobile view / GitHub problem" reflections as new paragraphs
# right after the paragraph that ends with "() to my image." (and before
# the trailing blank paragraph at the end of the document body).

$d = $word.ActiveDocument

$newParagraphs = @(
    "The mobile view is very basic",
    "For the mobile view this has the navigation in a list going down and also has four",
    "main content blocks underneath (two by two) so the user can scroll down and see",
    "past work rather than seeing one image at a time. This sketch has more padding around the main image and now has the navigation bar right at the top of the website as this can be viewed easily from the user. It also allows the main viewing point to be the image in the centre of the screen and the logo is top left. Following on from the last sketch this also has two images underneath rather than three as there will be writing at the top of them.",
    "In the begining of the project development process, I had to rename the project in GitHub unfortunately it overwrited my document with a blank one. I couldn't find solutions to recover it.",
    "Also in midway when I tried pushing to GitHub, I had error `"git.log is not recognized`" in my system. So I inquired this problem to Dipak sir, he gave me a solution to re-clone the existing repo. Before this, in regards to my previous problem, I had backed up my new file if the new file had been overwritten my the previous file from GitHub. As expected, thus I copied and pushed my back ups to the GitHub."
)

# Locate the paragraph that ends the "Additional:" reflection section --
# the one whose text contains "() to my image."
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*() to my image.*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph ending in '() to my image.'"
}

$currentIndex = $anchorIndex
foreach ($text in $newParagraphs) {
    $current = $d.Paragraphs.Item($currentIndex)
    $current.Range.InsertParagraphAfter()
    $currentIndex = $currentIndex + 1
    $newPara = $d.Paragraphs.Item($currentIndex)
    $newPara.Range.Text = $text
}

Write-Output ("Inserted " + $newParagraphs.Count + " paragraphs after paragraph " + $anchorIndex)
